$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# Insert 3 new columns (D:F) before the existing "Title" column, shifting it to G
$ws.Range("D1:F2").Insert(-4161)

# New header row (row 1)
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

# New data row (row 2)
$ws.Range("D2").Value = "Q (%)"
$ws.Range("E2").Value = "A (%)"
$ws.Range("F2").Value = "P (%)"

# Match the column widths applied to the new columns
$ws.Columns("D").ColumnWidth = 13.3
$ws.Columns("E").ColumnWidth = 13.3
$ws.Columns("F").ColumnWidth = 13.3

# Make "axes" the active sheet and select cell F8, matching the saved view state
$ws.Activate()
$ws.Range("F8").Select()
